$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 for Maidstone (pushes old row 8 "Melbourne" onward down by 1)
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Maidstone"
$ws.Range("B8").Value = "Marciano's Cakes  126 Mitchell St  Maidstone VIC 3012"
$ws.Range("C8").Value = "09:45-10:25  5/2/2021"
$ws.Range("D8").Value = "Case attended venue"

# After the insert above, the old rows 8-14 are now at rows 9-15.
# Insert 3 new rows after row 15 (Springvale / Woolworths) for Sunshine x2 and Taylors Lakes,
# pushing "West Melbourne" (now at row 15) down to row 19.
$ws.Range("16:18").Insert()

$ws.Range("A16").Value = "Sunshine"
$ws.Range("B16").Value = "Dan Murphy's  47 McIntyre Rd  Sunshine VIC 3020"
$ws.Range("C16").Value = "17:50-18:30  5/2/2021"
$ws.Range("D16").Value = "Case attended venue"

$ws.Range("A17").Value = "Sunshine"
$ws.Range("B17").Value = "Dan Murphy's  47 McIntyre Rd  Sunshine VIC 3020"
$ws.Range("C17").Value = "18:50-19:30  6/2/2021"
$ws.Range("D17").Value = "Case attended venue"

$ws.Range("A18").Value = "Taylors Lakes"
$ws.Range("B18").Value = "Off Ya Tree Watergardens  399 Melton Highway  Taylors Lakes VIC 3038"
$ws.Range("C18").Value = "13:17-13:52  6/2/2021"
$ws.Range("D18").Value = "Case attended venue"

$wb.Save()
